# Auto-generated Excel COM-interop script
# Applies the 2025-10-30 crime data update across all affected worksheets.
$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet1.xml)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 5580
$ws.Range("L3").Value = 6062
$ws.Range("K5").Value = 7653
$ws.Range("L5").Value = 5206
$ws.Range("E6").Value = 2062
$ws.Range("L6").Value = 1493
$ws.Range("L8").Value = 14207
$ws.Range("L9").Value = 4973
$ws.Range("I10").Value = 54901
$ws.Range("K10").Value = 61293
$ws.Range("L10").Value = 48749
$ws.Range("E11").Value = 113462
$ws.Range("I11").Value = 110708
$ws.Range("L11").Value = 86940

# By Neighborhood (sheet2.xml)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 763
$ws.Range("L4").Value = 515
$ws.Range("L7").Value = 2031
$ws.Range("L8").Value = 3762
$ws.Range("L11").Value = 1482
$ws.Range("L12").Value = 342
$ws.Range("L13").Value = 228
$ws.Range("L14").Value = 618
$ws.Range("L15").Value = 669
$ws.Range("L16").Value = 630
$ws.Range("L17").Value = 115
$ws.Range("L19").Value = 1979
$ws.Range("L20").Value = 1434
$ws.Range("L23").Value = 1108
$ws.Range("L25").Value = 495
$ws.Range("L27").Value = 1084
$ws.Range("L29").Value = 2842
$ws.Range("L30").Value = 225
$ws.Range("L31").Value = 835
$ws.Range("L32").Value = 152
$ws.Range("L33").Value = 2054
$ws.Range("L34").Value = 599
$ws.Range("L36").Value = 1199
$ws.Range("L37").Value = 1939
$ws.Range("L39").Value = 151
$ws.Range("L41").Value = 291
$ws.Range("L42").Value = 2096
$ws.Range("L43").Value = 1196
$ws.Range("L44").Value = 825
$ws.Range("L47").Value = 780
$ws.Range("L48").Value = 2311
$ws.Range("L49").Value = 1463
$ws.Range("L50").Value = 872
$ws.Range("L51").Value = 1226
$ws.Range("L52").Value = 1158
$ws.Range("L53").Value = 1471
$ws.Range("L54").Value = 3212
$ws.Range("L55").Value = 885
$ws.Range("L56").Value = 493
$ws.Range("L57").Value = 379
$ws.Range("L58").Value = 113
$ws.Range("L61").Value = 118
$ws.Range("E63").Value = 2113
$ws.Range("I63").Value = 1973
$ws.Range("L63").Value = 366
$ws.Range("L64").Value = 791
$ws.Range("L65").Value = 1161
$ws.Range("L66").Value = 497
$ws.Range("L67").Value = 1577
$ws.Range("L68").Value = 411
$ws.Range("L70").Value = 564
$ws.Range("L71").Value = 299
$ws.Range("L72").Value = 574
$ws.Range("L73").Value = 1119
$ws.Range("L74").Value = 266
$ws.Range("L75").Value = 352
$ws.Range("L76").Value = 2735
$ws.Range("L78").Value = 1321
$ws.Range("L79").Value = 1754
$ws.Range("L80").Value = 397
$ws.Range("L81").Value = 157
$ws.Range("L84").Value = 706
$ws.Range("L85").Value = 2838
$ws.Range("L86").Value = 776
$ws.Range("L87").Value = 374
$ws.Range("L88").Value = 726
$ws.Range("L89").Value = 1743
$ws.Range("L90").Value = 908
$ws.Range("L91").Value = 803
$ws.Range("L92").Value = 270
$ws.Range("L93").Value = 672
$ws.Range("L94").Value = 2248
$ws.Range("L95").Value = 1021
$ws.Range("L96").Value = 1133
$ws.Range("L97").Value = 1495
$ws.Range("L98").Value = 924
$ws.Range("L99").Value = 1250
$ws.Range("E101").Value = 113462
$ws.Range("I101").Value = 110708
$ws.Range("L101").Value = 86940

# Bridgeport (sheet3.xml)
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L5").Value = 57
$ws.Range("L10").Value = 356
$ws.Range("L11").Value = 618

# West Ridge (sheet4.xml)
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L3").Value = 60
$ws.Range("L8").Value = 208
$ws.Range("L10").Value = 628
$ws.Range("L11").Value = 1133

# Auburn Gresham (sheet5.xml)
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L5").Value = 163
$ws.Range("L10").Value = 776
$ws.Range("L11").Value = 2031

# Belmont Cragin (sheet6.xml)
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 117
$ws.Range("L8").Value = 219
$ws.Range("L10").Value = 878
$ws.Range("L11").Value = 1482

# O'Hare (sheet7.xml)
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("L3").Value = 15
$ws.Range("L11").Value = 564

# Uptown (sheet8.xml)
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L9").Value = 69
$ws.Range("L10").Value = 1207
$ws.Range("L11").Value = 1743

# South Shore (sheet9.xml)
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L8").Value = 531
$ws.Range("L10").Value = 1129
$ws.Range("L11").Value = 2838

# Little Village (sheet10.xml)
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L6").Value = 26
$ws.Range("L8").Value = 140
$ws.Range("L10").Value = 566
$ws.Range("L11").Value = 1158

# Logan Square (sheet12.xml)
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L8").Value = 187
$ws.Range("L10").Value = 994
$ws.Range("L11").Value = 1471

# Austin (sheet13.xml)
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 359
$ws.Range("L3").Value = 432
$ws.Range("L8").Value = 623
$ws.Range("L10").Value = 1704
$ws.Range("L11").Value = 3762

# Oakland (sheet16.xml)
$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("L8").Value = 101
$ws.Range("L11").Value = 299

# Garfield Park (sheet18.xml)
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L5").Value = 110
$ws.Range("L8").Value = 381
$ws.Range("L10").Value = 710
$ws.Range("L11").Value = 2054

# Roseland (sheet19.xml)
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 164
$ws.Range("L5").Value = 125
$ws.Range("L8").Value = 388
$ws.Range("L9").Value = 127
$ws.Range("L10").Value = 732
$ws.Range("L11").Value = 1754

# Pullman (sheet20.xml)
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("L8").Value = 54
$ws.Range("L11").Value = 352

# Near South Side (sheet22.xml)
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L8").Value = 190
$ws.Range("L11").Value = 791

# West Pullman (sheet23.xml)
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L5").Value = 118
$ws.Range("L8").Value = 243
$ws.Range("L10").Value = 386
$ws.Range("L11").Value = 1021

# Grand Crossing (sheet24.xml)
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 212
$ws.Range("L3").Value = 247
$ws.Range("L8").Value = 349
$ws.Range("L10").Value = 740
$ws.Range("L11").Value = 1939

# New City (sheet25.xml)
$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 131
$ws.Range("L5").Value = 94
$ws.Range("L11").Value = 1161

# Edgewater (sheet26.xml)
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L10").Value = 713
$ws.Range("L11").Value = 1084

# Woodlawn (sheet27.xml)
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 92
$ws.Range("L3").Value = 131
$ws.Range("L5").Value = 77
$ws.Range("L8").Value = 235
$ws.Range("L10").Value = 611
$ws.Range("L11").Value = 1250

# Fuller Park (sheet28.xml)
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L3").Value = 19
$ws.Range("L11").Value = 225

# Gage Park (sheet29.xml)
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L8").Value = 120
$ws.Range("L11").Value = 835

# North Lawndale (sheet30.xml)
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L8").Value = 260
$ws.Range("L9").Value = 147
$ws.Range("L10").Value = 595
$ws.Range("L11").Value = 1577

# South Deering (sheet31.xml)
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 57
$ws.Range("L8").Value = 160
$ws.Range("L10").Value = 337
$ws.Range("L11").Value = 706

# West Loop (sheet32.xml)
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L8").Value = 298
$ws.Range("L10").Value = 1641
$ws.Range("L11").Value = 2248

# River North (sheet33.xml)
$ws = $wb.Worksheets.Item('River North')
$ws.Range("L2").Value = 58
$ws.Range("L5").Value = 63
$ws.Range("L10").Value = 2215
$ws.Range("L11").Value = 2735

# Ukrainian Village (sheet34.xml)
$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("L10").Value = 254
$ws.Range("L11").Value = 374

# East Side (sheet35.xml)
$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L10").Value = 182
$ws.Range("L11").Value = 495

# Bucktown (sheet37.xml)
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("L10").Value = 494
$ws.Range("L11").Value = 630

# Lincoln Park (sheet38.xml)
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L8").Value = 130
$ws.Range("L10").Value = 1138
$ws.Range("L11").Value = 1463

# West Town (sheet39.xml)
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L6").Value = 12
$ws.Range("L8").Value = 269
$ws.Range("L10").Value = 976
$ws.Range("L11").Value = 1495

# Lower West Side (sheet40.xml)
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L5").Value = 63
$ws.Range("L10").Value = 470
$ws.Range("L11").Value = 885

# Loop (sheet41.xml)
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L8").Value = 303
$ws.Range("L10").Value = 2435
$ws.Range("L11").Value = 3212

# Portage Park (sheet42.xml)
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L6").Value = 16
$ws.Range("L10").Value = 742
$ws.Range("L11").Value = 1119

# Englewood (sheet43.xml)
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L5").Value = 184
$ws.Range("L6").Value = 55
$ws.Range("L8").Value = 545
$ws.Range("L9").Value = 261
$ws.Range("L10").Value = 1057
$ws.Range("L11").Value = 2842

# Lake View (sheet44.xml)
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L8").Value = 205
$ws.Range("L10").Value = 1674
$ws.Range("L11").Value = 2311

# Chatham (sheet45.xml)
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L5").Value = 156
$ws.Range("L10").Value = 976
$ws.Range("L11").Value = 1979

# North Center (sheet46.xml)
$ws = $wb.Worksheets.Item('North Center')
$ws.Range("L10").Value = 317
$ws.Range("L11").Value = 497

# Irving Park (sheet47.xml)
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L8").Value = 127
$ws.Range("L10").Value = 515
$ws.Range("L11").Value = 825

# Humboldt Park (sheet48.xml)
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 162
$ws.Range("L3").Value = 207
$ws.Range("L8").Value = 396
$ws.Range("L10").Value = 989
$ws.Range("L11").Value = 2096

# Hermosa (sheet51.xml)
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L8").Value = 52
$ws.Range("L10").Value = 137
$ws.Range("L11").Value = 291

# Grand Boulevard (sheet52.xml)
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L5").Value = 63
$ws.Range("L8").Value = 371
$ws.Range("L10").Value = 529
$ws.Range("L11").Value = 1199

# Boystown (sheet53.xml)
$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("L2").Value = 3
$ws.Range("L9").Value = 175
$ws.Range("L10").Value = 228

# Streeterville (sheet55.xml)
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L6").Value = 68
$ws.Range("L8").Value = 79
$ws.Range("L11").Value = 776

# Rogers Park (sheet56.xml)
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 79
$ws.Range("L10").Value = 817
$ws.Range("L11").Value = 1321

# North Park (sheet57.xml)
$ws = $wb.Worksheets.Item('North Park')
$ws.Range("L6").Value = 7
$ws.Range("L8").Value = 43
$ws.Range("L10").Value = 276
$ws.Range("L11").Value = 411

# Brighton Park (sheet59.xml)
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L8").Value = 99
$ws.Range("L10").Value = 363
$ws.Range("L11").Value = 669

# Douglas (sheet60.xml)
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K5").Value = 50
$ws.Range("L6").Value = 17
$ws.Range("L8").Value = 283
$ws.Range("K10").Value = 726
$ws.Range("L11").Value = 1108

# Washington Park (sheet61.xml)
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L10").Value = 303
$ws.Range("L11").Value = 803

# Little Italy, UIC (sheet62.xml)
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L8").Value = 251
$ws.Range("L10").Value = 692
$ws.Range("L11").Value = 1226

# Chicago Lawn (sheet63.xml)
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L8").Value = 228
$ws.Range("L10").Value = 631
$ws.Range("L11").Value = 1434

# Kenwood (sheet65.xml)
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L3").Value = 40
$ws.Range("L8").Value = 165
$ws.Range("L10").Value = 457
$ws.Range("L11").Value = 780

# Washington Heights (sheet66.xml)
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L5").Value = 72
$ws.Range("L8").Value = 240
$ws.Range("L10").Value = 403
$ws.Range("L11").Value = 908

# Lincoln Square (sheet67.xml)
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L10").Value = 634
$ws.Range("L11").Value = 872

# West Lawn (sheet68.xml)
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L8").Value = 129
$ws.Range("L10").Value = 402
$ws.Range("L11").Value = 672

# Magnificent Mile (sheet71.xml)
$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("L5").Value = 6
$ws.Range("L10").Value = 466
$ws.Range("L11").Value = 493

# Albany Park (sheet73.xml)
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L3").Value = 51
$ws.Range("L5").Value = 38
$ws.Range("L8").Value = 131
$ws.Range("L9").Value = 40
$ws.Range("L10").Value = 439
$ws.Range("L11").Value = 763

# Old Town (sheet74.xml)
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L10").Value = 386
$ws.Range("L11").Value = 574

# Hyde Park (sheet75.xml)
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L2").Value = 28
$ws.Range("L10").Value = 855
$ws.Range("L11").Value = 1196

# Burnside (sheet76.xml)
$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("L10").Value = 48
$ws.Range("L11").Value = 115

# Archer Heights (sheet77.xml)
$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L10").Value = 348
$ws.Range("L11").Value = 515

# Rush & Division (sheet78.xml)
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("L10").Value = 314
$ws.Range("L11").Value = 397

# Garfield Ridge (sheet79.xml)
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L10").Value = 328
$ws.Range("L11").Value = 599

# Galewood (sheet82.xml)
$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("L10").Value = 80
$ws.Range("L11").Value = 152

# Millenium Park (sheet83.xml)
$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("L10").Value = 101
$ws.Range("L11").Value = 113

# Wicker Park (sheet84.xml)
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L8").Value = 132
$ws.Range("L10").Value = 639
$ws.Range("L11").Value = 924

# West Elsdon (sheet86.xml)
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L10").Value = 132
$ws.Range("L11").Value = 270

# Mount Greenwood (sheet87.xml)
$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("L5").Value = 12
$ws.Range("L11").Value = 118

# Greektown (sheet91.xml)
$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("L8").Value = 122
$ws.Range("L9").Value = 151

# Mckinley Park (sheet92.xml)
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("L3").Value = 20
$ws.Range("L8").Value = 48
$ws.Range("L10").Value = 237
$ws.Range("L11").Value = 379

# United Center (sheet93.xml)
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 68
$ws.Range("L10").Value = 355
$ws.Range("L11").Value = 726

# Printers Row (sheet94.xml)
$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("L3").Value = 4
$ws.Range("L11").Value = 266

# Sauganash,Forest Glen (sheet96.xml)
$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("L10").Value = 95
$ws.Range("L11").Value = 157

# Beverly (sheet98.xml)
$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L10").Value = 213
$ws.Range("L11").Value = 342

